$wb = $excel.ActiveWorkbook

# Rename sheets for parallel structure with the data-processing R script
$wb.Worksheets.Item("1881").Name = "1880Survey"
$wb.Worksheets.Item("1940").Name = "1940Survey"

$ws1 = $wb.Worksheets.Item("1880Survey")
$ws2 = $wb.Worksheets.Item("1881notes")
$ws3 = $wb.Worksheets.Item("1940Survey")
$ws4 = $wb.Worksheets.Item("1940notes")

# Add a new metadata / ToDo sheet at the very end of the workbook
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws4)
$newSheet.Name = "1880Metadata"
$newSheet.Range("A1").Value = "Actually surveyed in 1881; tab reads 1880 for consistency in the data processing R script"

# Fill Survey_id (column B) for rows 3-40 on the 1880Survey sheet to match row 2
$surveyId = $ws1.Range("B2").Value2
for ($r = 3; $r -le 40; $r++) {
    $ws1.Cells.Item($r, 2).Value = $surveyId
}

# Restore each sheet's selection state
$ws2.Activate() | Out-Null
$ws2.Range("A8").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("G30").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("A3").Select() | Out-Null

$newSheet.Activate() | Out-Null
$newSheet.Range("A2").Select() | Out-Null

# 1880Survey is the tab that should be selected/active when the workbook opens
$ws1.Activate() | Out-Null
$ws1.Range("C1").Select() | Out-Null
